$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.262.70'
$ws.Range("E2").Value = '  -0.11%  '
$ws.Range("D3").Value = '1.684.23'
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.89'
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5251'
$ws.Range("E6").Value = '  +2.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.008'
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2704'
$ws.Range("E8").Value = '  +1.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06413'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.97'
$ws.Range("E10").Value = '  +1.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07487'
$ws.Range("E11").Value = '  +1.54%  '
$ws.Range("D12").Value = '1.722.18'
$ws.Range("E12").Value = '  +2.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.549'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5803'
$ws.Range("E14").Value = '  -0.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008465'
$ws.Range("E15").Value = '  -2.69%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.19'
$ws.Range("E16").Value = '  -1.30%  '
$ws.Range("D17").Value = '26.323.57'
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.923'
$ws.Range("E18").Value = '  -1.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.007'
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.85'
$ws.Range("E20").Value = '  -0.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '188.99'
$ws.Range("E21").Value = '  -1.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.194'
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.36'
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.706'
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1235'
$ws.Range("E26").Value = '  +4.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.77'
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06656'
$ws.Range("E28").Value = '  +12.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.349'
$ws.Range("E29").Value = '  +5.84%  '
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.570'
$ws.Range("E31").Value = '  +1.29%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.564'
$ws.Range("E32").Value = '  +0.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.663'
$ws.Range("E33").Value = '  +0.84%  '
$ws.Range("E34").Value = '  +1.19%  '
$ws.Range("E35").Value = '  +2.80%  '
$ws.Range("E36").Value = '  +1.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.704'
$ws.Range("E37").Value = '  +2.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.390'
$ws.Range("E38").Value = '  +5.63%  '
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '1.104.84'
$ws.Range("E40").Value = '  +2.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8764'
$ws.Range("E41").Value = '  +0.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.015'
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.62'
$ws.Range("E43").Value = '  +0.76%  '
$ws.Range("D44").Value = '1.832.35'
$ws.Range("E44").Value = '  +0.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("E45").Value = '  -2.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.70'
$ws.Range("E46").Value = '  +0.91%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.009'
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.153'
$ws.Range("E48").Value = '  +0.43%  '
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4304'
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.027'
$ws.Range("E51").Value = '  +2.19%  '
